$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 9
$wsALC.Range("H9").Value = 333700.34
$wsALC.Range("I9").Value = 550
$wsALC.Range("J9").Value = 1000001
$wsALC.Range("K9").Value = 550
$wsALC.Range("L9").Value = 1000001
$wsALC.Range("M9").Value = -381
$wsALC.Range("N9").Value = -1000339

# ALC row 29
$wsALC.Range("H29").Value = 313.25
$wsALC.Range("I29").Value = 313.25
$wsALC.Range("K29").Value = 939.75
$wsALC.Range("M29").Value = -658.75

# ALC row 38
$wsALC.Range("H38").Value = 302.66666
$wsALC.Range("I38").Value = 183.2
$wsALC.Range("K38").Value = 549.5999999999999
$wsALC.Range("M38").Value = -177.5999999999999

# ALC row 58
$wsALC.Range("H58").Value = 1666.2174
$wsALC.Range("I58").Value = 1174.8182
$wsALC.Range("J58").Value = 2116.6667
$wsALC.Range("K58").Value = 3524.4546
$wsALC.Range("L58").Value = 6350.000100000001
$wsALC.Range("M58").Value = -3374.4546
$wsALC.Range("N58").Value = -6650.000100000001

# ALC row 80
$wsALC.Range("H80").Value = 0
$wsALC.Range("I80").Value = 0
$wsALC.Range("J80").Value = 0
$wsALC.Range("K80").Value = 0
$wsALC.Range("L80").Value = 0
$wsALC.Range("M80").ClearContents()
$wsALC.Range("N80").ClearContents()

# ALC row 83
$wsALC.Range("H83").Value = 0
$wsALC.Range("I83").Value = 0
$wsALC.Range("J83").Value = 0
$wsALC.Range("K83").Value = 0
$wsALC.Range("L83").Value = 0
$wsALC.Range("M83").ClearContents()
$wsALC.Range("N83").ClearContents()

# ALC row 87
$wsALC.Range("H87").Value = 23999.436
$wsALC.Range("J87").Value = 23999.436
$wsALC.Range("L87").Value = 23999.436
$wsALC.Range("N87").Value = -26495.436

# ALC row 90
$wsALC.Range("H90").Value = 23999.436
$wsALC.Range("J90").Value = 23999.436
$wsALC.Range("L90").Value = 71998.308
$wsALC.Range("N90").Value = -84478.308

# ALC row 132
$wsALC.Range("H132").Value = 2583.8108
$wsALC.Range("I132").Value = 2557.9395
$wsALC.Range("J132").Value = 2797.25
$wsALC.Range("K132").Value = 7673.818499999999
$wsALC.Range("L132").Value = 8391.75
$wsALC.Range("M132").Value = -5143.818499999999
$wsALC.Range("N132").Value = -13451.75

# ALC row 138
$wsALC.Range("H138").Value = 2101.82
$wsALC.Range("J138").Value = 2914.1052
$wsALC.Range("L138").Value = 8742.3156
$wsALC.Range("N138").Value = -19022.3156

# ARM row 32
$wsARM.Range("H32").Value = 488870.2
$wsARM.Range("I32").Value = 537046
$wsARM.Range("K32").Value = 537046
$wsARM.Range("M32").Value = -536759

# ARM row 61
$wsARM.Range("H61").Value = 2293.5264
$wsARM.Range("I61").Value = 1946.9032
$wsARM.Range("K61").Value = 1946.9032
$wsARM.Range("M61").Value = -1734.9032

# ARM row 74
$wsARM.Range("H74").Value = 1054.7567
$wsARM.Range("I74").Value = 736.0357
$wsARM.Range("J74").Value = 2046.3334
$wsARM.Range("K74").Value = 736.0357
$wsARM.Range("L74").Value = 2046.3334
$wsARM.Range("M74").Value = 137.9643
$wsARM.Range("N74").Value = -3794.3334

# ARM row 77
$wsARM.Range("H77").Value = 1054.7567
$wsARM.Range("I77").Value = 736.0357
$wsARM.Range("J77").Value = 2046.3334
$wsARM.Range("K77").Value = 3680.1785
$wsARM.Range("L77").Value = 10231.667
$wsARM.Range("M77").Value = 687.8215
$wsARM.Range("N77").Value = -18967.667

# ARM row 122
$wsARM.Range("H122").Value = 78255.69500000001
$wsARM.Range("I122").Value = 100882.4
$wsARM.Range("K122").Value = 302647.2
$wsARM.Range("M122").Value = -300197.2

# ARM row 128
$wsARM.Range("H128").Value = 39909.668
$wsARM.Range("J128").Value = 39909.668
$wsARM.Range("L128").Value = 39909.668
$wsARM.Range("N128").Value = -49869.668

# ARM row 132
$wsARM.Range("H132").Value = 4138.8335
$wsARM.Range("I132").Value = 2686.9565
$wsARM.Range("J132").Value = 6707.5386
$wsARM.Range("K132").Value = 8060.869499999999
$wsARM.Range("L132").Value = 20122.6158
$wsARM.Range("M132").Value = -5530.869499999999
$wsARM.Range("N132").Value = -25182.6158

# ARM row 136
$wsARM.Range("H136").Value = 2293.5264
$wsARM.Range("I136").Value = 1946.9032
$wsARM.Range("K136").Value = 5840.7096
$wsARM.Range("M136").Value = -3290.7096

# BSM row 74
$wsBSM.Range("H74").Value = 20000
$wsBSM.Range("J74").Value = 20000
$wsBSM.Range("L74").Value = 20000
$wsBSM.Range("N74").Value = -21872

# BSM row 77
$wsBSM.Range("H77").Value = 20000
$wsBSM.Range("J77").Value = 20000
$wsBSM.Range("L77").Value = 60000
$wsBSM.Range("N77").Value = -69360

# BSM row 99
$wsBSM.Range("H99").Value = 1842
$wsBSM.Range("I99").Value = 1759.3
$wsBSM.Range("J99").Value = 2255.5
$wsBSM.Range("K99").Value = 1759.3
$wsBSM.Range("L99").Value = 2255.5
$wsBSM.Range("M99").Value = -261.3
$wsBSM.Range("N99").Value = -5251.5

# BSM row 134
$wsBSM.Range("H134").Value = 2204.6667
$wsBSM.Range("I134").Value = 1477.1765
$wsBSM.Range("J134").Value = 3971.4285
$wsBSM.Range("K134").Value = 4431.529500000001
$wsBSM.Range("L134").Value = 11914.2855
$wsBSM.Range("M134").Value = -1896.529500000001
$wsBSM.Range("N134").Value = -16984.2855

# BSM row 141
$wsBSM.Range("H141").Value = 24998.334
$wsBSM.Range("J141").Value = 24998.334
$wsBSM.Range("L141").Value = 24998.334
$wsBSM.Range("N141").Value = -35358.334

# CRP row 62
$wsCRP.Range("H62").Value = 1833.3334
$wsCRP.Range("I62").Value = 1000
$wsCRP.Range("J62").Value = 2250
$wsCRP.Range("K62").Value = 1000
$wsCRP.Range("L62").Value = 2250
$wsCRP.Range("M62").Value = -376
$wsCRP.Range("N62").Value = -3498

# CRP row 65
$wsCRP.Range("H65").Value = 1833.3334
$wsCRP.Range("I65").Value = 1000
$wsCRP.Range("J65").Value = 2250
$wsCRP.Range("K65").Value = 5000
$wsCRP.Range("L65").Value = 11250
$wsCRP.Range("M65").Value = -1880
$wsCRP.Range("N65").Value = -17490

# CRP row 122
$wsCRP.Range("H122").Value = 1902.8718
$wsCRP.Range("I122").Value = 1821.1
$wsCRP.Range("J122").Value = 1988.9474
$wsCRP.Range("K122").Value = 5463.299999999999
$wsCRP.Range("L122").Value = 5966.8422
$wsCRP.Range("M122").Value = -3013.299999999999
$wsCRP.Range("N122").Value = -10866.8422

# CUL row 12
$wsCUL.Range("H12").Value = 115.82353
$wsCUL.Range("I12").Value = 71.333336
$wsCUL.Range("K12").Value = 214.000008
$wsCUL.Range("M12").Value = -41.00000800000001

# CUL row 23
$wsCUL.Range("H23").Value = 88235430
$wsCUL.Range("I23").Value = 137.6
$wsCUL.Range("J23").Value = 125000140
$wsCUL.Range("K23").Value = 412.8
$wsCUL.Range("L23").Value = 375000420
$wsCUL.Range("M23").Value = -177.8
$wsCUL.Range("N23").Value = -375000890

# CUL row 31
$wsCUL.Range("H31").Value = 1417.091
$wsCUL.Range("J31").Value = 1258.8
$wsCUL.Range("L31").Value = 3776.4
$wsCUL.Range("N31").Value = -4352.4

# CUL row 132
$wsCUL.Range("H132").Value = 2733.0417
$wsCUL.Range("J132").Value = 2713.2778
$wsCUL.Range("L132").Value = 24419.5002
$wsCUL.Range("N132").Value = -29479.5002

# WVR row 81
$wsWVR.Range("H81").Value = 6115.273
$wsWVR.Range("J81").Value = 4524.75
$wsWVR.Range("L81").Value = 9049.5
$wsWVR.Range("N81").Value = -11171.5

# WVR row 84
$wsWVR.Range("H84").Value = 6115.273
$wsWVR.Range("J84").Value = 4524.75
$wsWVR.Range("L84").Value = 45247.5
$wsWVR.Range("N84").Value = -55855.5

# WVR row 132
$wsWVR.Range("H132").Value = 1517.1025
$wsWVR.Range("I132").Value = 1355.2174
$wsWVR.Range("J132").Value = 1749.8125
$wsWVR.Range("K132").Value = 4065.6522
$wsWVR.Range("L132").Value = 5249.4375
$wsWVR.Range("M132").Value = -1535.6522
$wsWVR.Range("N132").Value = -10309.4375
